# Template fix: collapse the "{ tag }" runs (which Word's spell-checker had
# split up with <w:proofErr/> wrappers around the bare tag name) back into a
# single run of literal text per placeholder, and label the first Likert
# table column "Question" instead of a lone non-breaking space.

$d = $word.ActiveDocument

# --- Simple placeholder merges -------------------------------------------
# Each of these previously spanned 3 runs (prefix, proofErr-wrapped tag
# name, suffix). Find/Replace re-merges the matched text into one run
# (clearing the spell-check proofErr markers in the process), which is
# exactly what the diff shows.

$d.Content.Find.Execute(": {childName}", $true, $false, $false, $false, $false, `
    $true, 1, $false, ": {childName}", 2) | Out-Null

$d.Content.Find.Execute(": {dateOfBirth}", $true, $false, $false, $false, $false, `
    $true, 1, $false, ": {dateOfBirth}", 2) | Out-Null

$d.Content.Find.Execute("{pregWeeks}", $true, $false, $false, $false, $false, `
    $true, 1, $false, "{pregWeeks}", 2) | Out-Null

$d.Content.Find.Execute("{birthProblems}", $true, $false, $false, $false, $false, `
    $true, 1, $false, "{birthProblems}", 2) | Out-Null

$d.Content.Find.Execute("{affectingConditions}", $true, $false, $false, $false, $false, `
    $true, 1, $false, "{affectingConditions}", 2) | Out-Null

$d.Content.Find.Execute("{visionConcerns}", $true, $false, $false, $false, $false, `
    $true, 1, $false, "{visionConcerns}", 2) | Out-Null

$d.Content.Find.Execute("{not_applicable} {/likertTable}", $true, $false, $false, $false, $false, `
    $true, 1, $false, "{not_applicable} {/likertTable}", 2) | Out-Null

# --- Likert table header: blank cell -> "Question" -----------------------
# The cell held a single non-breaking space (U+00A0) run. Replace it with
# the full word "Question" (Find keeps the existing bold run formatting),
# then split that single run into "Questio" + "n" by nudging the Bold
# property of the last character off and back on -- two adjacent runs with
# identical formatting get coalesced by the engine unless their formatting
# genuinely differs at some point, so toggling forces the split while the
# final state matches the original (bold) formatting exactly.

$found = $d.Content.Find.Execute([char]0x00A0, $true, $false, $false, $false, $false, `
    $true, 1, $false, "Question", 2)

if ($found) {
    $headerPara = $d.Paragraphs.Item(10)
    $headerRange = $headerPara.Range
    $lastCharStart = $headerRange.Start + 7
    $lastCharEnd = $lastCharStart + 1

    $lastChar = $d.Range($lastCharStart, $lastCharEnd)
    $lastChar.Bold = 0

    $lastChar2 = $d.Range($lastCharStart, $lastCharEnd)
    $lastChar2.Bold = 1
}
